$wb = $excel.ActiveWorkbook

$newAddress = "括苍路493号油泵厂山顶通用设备厂区块3号楼 中国国际摄影节展览馆"
$newCount = 1011

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D2").Value = $newAddress
    $ws.Range("F2").Value = $newCount
}
